$d = $word.ActiveDocument

$replacements = @(
    @("320×6=", "717×2="),
    @("234×2=", "285×9="),
    @("884×3=", "845×7="),
    @("572×8=", "643×5="),
    @("164×7=", "888×6="),
    @("224×7=", "889×2="),
    @("301×3=", "410×3="),
    @("512×7=", "962×7="),
    @("221×8=", "734×5="),
    @("348×7=", "589×2="),
    @("847×5=", "132×8="),
    @("812×4=", "450×6="),
    @("509×8=", "305×2="),
    @("312×4=", "657×7="),
    @("316×8=", "977×6="),
    @("552×6=", "979×9="),
    @("855×4=", "251×3="),
    @("507×4=", "884×6="),
    @("869×9=", "202×6="),
    @("640×4=", "275×5="),
    @("371×2=", "466×3="),
    @("307×9=", "839×9="),
    @("264×7=", "923×6="),
    @("842×7=", "247×7="),
    @("578×8=", "831×7=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
